$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("I1").Value = "Size"
$ws.Range("I1").Font.Bold = $true

$sizes = @("L", "M", "XL", "S", "M", "L", "XS", "M", "XL")
for ($i = 0; $i -lt $sizes.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 9).Value = $sizes[$i]
}

$ws.Range("I10").Select()
